# feat: add 2022-Q4 data
#
# Before:  Sheet1 "总计" (totals) | Sheet2 "2022-Q1" (fund-holder detail)
# After:   Sheet1 "总计" (totals, now with a Q4 row too)
#          Sheet2 "2022-Q4" (new fund-holder detail for Q4)
#          Sheet3 "2022-Q1" (the original fund-holder detail, unchanged, relocated)

$wb  = $excel.ActiveWorkbook
$tot = $wb.Worksheets.Item(1)
$oldQ1 = $wb.Worksheets.Item(2)

# ------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q1" detail sheet *before* touching it,
#    so the duplicate preserves the original data/formatting untouched.
#    The duplicate lands right after the original and keeps tab-selection.
# ------------------------------------------------------------------
$oldQ1.Copy($null, $oldQ1)
$q4   = $wb.Worksheets.Item(2)   # the original physical sheet -> becomes 2022-Q4
$newQ1 = $wb.Worksheets.Item(3)  # the fresh duplicate -> stays as 2022-Q1

$q4.Name = "2022-Q4"
$newQ1.Name = "2022-Q1"

# ------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: row 2 switches from 2022-Q1 to
#    2022-Q4 figures, and a new row 3 restores the original 2022-Q1 totals.
# ------------------------------------------------------------------
$tot.Range("B2").Value = "2022-Q4"
$tot.Range("C2").Value = 2
$tot.Range("D2").Value = 0.04

$tot.Range("A3").Value = 1
$tot.Range("B3").Value = "2022-Q1"
$tot.Range("C3").Value = 2
$tot.Range("D3").Value = 0.2

# Match formatting of the new row to the existing row above it.
$tot.Range("A2").Copy()
$tot.Range("A3").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3) Replace the "2022-Q4" sheet's contents with the new quarter's
#    fund-holder detail (it currently still holds the old Q1 rows).
# ------------------------------------------------------------------
$q4.Range("A1:H3").ClearContents()

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'016307"
$q4.Range("C2").Value = "景顺长城北交所精选两年定开混合A"
$q4.Range("D2").Value = "'1.83"
$q4.Range("E2").Value = "'43.56"
$q4.Range("F2").Value = "'2.04"
$q4.Range("G2").Value = "'0.0373"
$q4.Range("H2").Value = 9

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'016308"
$q4.Range("C3").Value = "景顺长城北交所精选两年定开混合C"
$q4.Range("D3").Value = "'0.27"
$q4.Range("E3").Value = "'43.56"
$q4.Range("F3").Value = "'2.04"
$q4.Range("G3").Value = "'0.0055"
$q4.Range("H3").Value = 9

# The new Q4 sheet is built from the "总计" template formatting (bold
# header/first-column style + its page margins), not the old Q1 sheet's.
$tot.Range("B1:D1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$tot.Range("A2").Copy()
$q4.Range("A2:A3").PasteSpecial(-4122)

$q4.PageSetup.LeftMargin = $tot.PageSetup.LeftMargin
$q4.PageSetup.RightMargin = $tot.PageSetup.RightMargin
$q4.PageSetup.TopMargin = $tot.PageSetup.TopMargin
$q4.PageSetup.BottomMargin = $tot.PageSetup.BottomMargin
$q4.PageSetup.HeaderMargin = $tot.PageSetup.HeaderMargin
$q4.PageSetup.FooterMargin = $tot.PageSetup.FooterMargin

# ------------------------------------------------------------------
# 4) Keep the relocated "2022-Q1" sheet as the active tab, matching the
#    original workbook (its tab was the selected one before the edit).
# ------------------------------------------------------------------
$newQ1.Activate()
$newQ1.Range("A1").Select()
